$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.770.67"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.90%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.606.15"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.35%  "

# Row 4
$ws.Range("E4").Value = "  +0.18%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "514.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.92%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.06"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.43%  "

# Row 7
$ws.Range("E7").Value = "  +0.23%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.563"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.43%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.610.30"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.18%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.26"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.70%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.103"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.22%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.337"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.39%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.128"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.89%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.072.78"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.73%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.807.99"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.84%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.94"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.44%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000136"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.41%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.619.01"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.00%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.38%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "339.85"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.67%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.84%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.04"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.00%  "

# Row 23
$ws.Range("E23").Value = "  -0.19%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.97"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.17%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.413"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.15%  "

# Row 26
$ws.Range("E26").Value = "  +0.34%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.157"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.23%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0793"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.89%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.92"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.03%  "

# Row 30
$ws.Range("E30").Value = "  +0.07%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.57"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.40%  "

# Row 32
$ws.Range("B32").Value = "Aptos"
$ws.Range("C32").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.90"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.44%  "

# Row 33
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.73"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.23%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "150.10"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.93%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.88"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.28%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.900"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.83%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.12"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.77%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.62"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.39%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.841"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.91%  "

# Row 40
$ws.Range("E40").Value = "  -6.24%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.57"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.16%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "283.22"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.38%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.622"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.29%  "

# Row 44
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.33%  "

# Row 45
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0990"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.43%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0539"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.80%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.02"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.92%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.39"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.80%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0229"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.53%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.58"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.06%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.936.30"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.06%  "

